# Update countries & provincias Spain
# Applies the 31-Mar-2020 22:50 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 22:50"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 183930
$ws.Range("C4").Value = 20142
$ws.Range("D4").Value = 6275
$ws.Range("E4").Value = 173932
$ws.Range("G4").Value = 582
$ws.Range("H4").Value = 3723

# Canada (row 18)
$ws.Range("B18").Value = 8505
$ws.Range("C18").Value = 1057
$ws.Range("E18").Value = 7248

# Row 20: Israel -> Brasil (with new, larger Brasil figures; Brasil overtakes Israel in ranking)
$ws.Range("A20").Value = "Brasil"
$ws.Range("B20").Value = 5717
$ws.Range("C20").Value = 1087
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 5389
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 201

# Row 21: Brasil -> Israel (keeps the former Israel figures, now re-ranked to row 21)
$ws.Range("A21").Value = "Israel"
$ws.Range("B21").Value = 5358
$ws.Range("C21").Value = 663
$ws.Range("D21").Value = 224
$ws.Range("E21").Value = 5114
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 20

# Sudafrica (row 43)
$ws.Range("E43").Value = 1317
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 5

# Oman (row 94)
$ws.Range("E94").Value = 157
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 1

# Estado de Palestina (row 111)
$ws.Range("B111").Value = 119
$ws.Range("C111").Value = 2
$ws.Range("E111").Value = 100

# Trinidad yTobago (row 120)
$ws.Range("B120").Value = 87
$ws.Range("C120").Value = 4
$ws.Range("E120").Value = 83
